# Update cryptos list (GitHub Actions scheduled refresh).
# Writes new Price (D) / Volume(1h) (E) values, and for rows 29/30 and
# 43/44 also refreshes Coin (B) and Link (C) since the two rows swapped
# rank order in the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to stay text even when the new value looks like a
    # number (e.g. '1.00', '537.22'); Range.Value would otherwise coerce
    # it to a numeric cell type and drop formatting like trailing zeros.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '59.136.08'
$ws.Range("E2").Value = '  +1.09%  '

# Row 3
$ws.Range("D3").Value = '2.522.25'
$ws.Range("E3").Value = '  +2.50%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
Set-TextValue "D5" '537.22'
$ws.Range("E5").Value = '  +1.18%  '

# Row 6
Set-TextValue "D6" '136.30'
$ws.Range("E6").Value = '  +1.03%  '

# Row 7
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  +0.19%  '

# Row 8
Set-TextValue "D8" '0.567'
$ws.Range("E8").Value = '  +1.74%  '

# Row 9
$ws.Range("D9").Value = '2.522.04'
$ws.Range("E9").Value = '  +2.10%  '

# Row 10
$ws.Range("E10").Value = '  +2.20%  '

# Row 11
$ws.Range("E11").Value = '  -2.14%  '

# Row 12
$ws.Range("E12").Value = '  -1.22%  '

# Row 13
$ws.Range("E13").Value = '  +0.88%  '

# Row 14
$ws.Range("D14").Value = '2.967.05'
$ws.Range("E14").Value = '  +2.35%  '

# Row 15
Set-TextValue "D15" '23.02'
$ws.Range("E15").Value = '  +1.89%  '

# Row 16
$ws.Range("D16").Value = '59.042.66'
$ws.Range("E16").Value = '  +1.09%  '

# Row 17
$ws.Range("E17").Value = '  +0.49%  '

# Row 18
$ws.Range("D18").Value = '2.521.03'
$ws.Range("E18").Value = '  +2.19%  '

# Row 19
Set-TextValue "D19" '11.08'
$ws.Range("E19").Value = '  +2.92%  '

# Row 20
Set-TextValue "D20" '4.27'
$ws.Range("E20").Value = '  +1.82%  '

# Row 21
Set-TextValue "D21" '322.54'
$ws.Range("E21").Value = '  +0.72%  '

# Row 22
Set-TextValue "D22" '1.00'
$ws.Range("E22").Value = '  +0.26%  '

# Row 23
Set-TextValue "D23" '5.97'
$ws.Range("E23").Value = '  +4.40%  '

# Row 24
Set-TextValue "D24" '65.06'
$ws.Range("E24").Value = '  +4.52%  '

# Row 25
$ws.Range("E25").Value = '  +3.14%  '

# Row 26
$ws.Range("E26").Value = '  -0.90%  '

# Row 27
Set-TextValue "D27" '0.998'
$ws.Range("E27").Value = '  +0.87%  '

# Row 28
Set-TextValue "D28" '7.52'
$ws.Range("E28").Value = '  -0.51%  '

# Row 29
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D29" '6.67'
$ws.Range("E29").Value = '  +2.20%  '

# Row 30
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0768'
$ws.Range("E30").Value = '  +1.37%  '

# Row 31
$ws.Range("E31").Value = '  -0.43%  '

# Row 32
Set-TextValue "D32" '170.58'
$ws.Range("E32").Value = '  +4.04%  '

# Row 33
$ws.Range("E33").Value = '  +8.73%  '

# Row 34
Set-TextValue "D34" '0.999'
$ws.Range("E34").Value = '  +0.03%  '

# Row 35
$ws.Range("E35").Value = '  +1.30%  '

# Row 36
Set-TextValue "D36" '18.36'
$ws.Range("E36").Value = '  +0.73%  '

# Row 37
$ws.Range("E37").Value = '  +0.74%  '

# Row 38
Set-TextValue "D38" '1.53'
$ws.Range("E38").Value = '  -0.68%  '

# Row 39
Set-TextValue "D39" '36.90'
$ws.Range("E39").Value = '  +1.15%  '

# Row 40
Set-TextValue "D40" '0.811'
$ws.Range("E40").Value = '  +3.43%  '

# Row 41
$ws.Range("E41").Value = '  +1.15%  '

# Row 42
Set-TextValue "D42" '284.92'
$ws.Range("E42").Value = '  +3.43%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D43" '0.999'
$ws.Range("E43").Value = '  +0.10%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D44" '5.08'
$ws.Range("E44").Value = '  -0.21%  '

# Row 45
$ws.Range("E45").Value = '  +3.24%  '

# Row 46
Set-TextValue "D46" '130.49'
$ws.Range("E46").Value = '  +7.15%  '

# Row 47
Set-TextValue "D47" '10.90'
$ws.Range("E47").Value = '  +0.54%  '

# Row 48
Set-TextValue "D48" '0.0924'
$ws.Range("E48").Value = '  -0.52%  '

# Row 49
$ws.Range("E49").Value = '  +0.17%  '

# Row 50
$ws.Range("E50").Value = '  -0.14%  '

# Row 51
Set-TextValue "D51" '17.42'
$ws.Range("E51").Value = '  +1.92%  '
